$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 16; this shifts existing rows 16-38 down to 17-39,
# carrying their formatting/values along with them.
$ws.Rows.Item(16).Insert()

# Populate the newly inserted row 16 with the new weekly record.
# (Same dimension values as the record that used to occupy row 16, but with an
# updated date and volume, matching the new weekly observation.)
$ws.Cells.Item(16, 1).Value = 7
$ws.Cells.Item(16, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(16, 3).Value = "Ñuble"
$ws.Cells.Item(16, 4).Value = (Get-Date -Year 2023 -Month 11 -Day 13 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(16, 4).Style = $ws.Cells.Item(17, 4).Style
$ws.Cells.Item(16, 4).NumberFormat = $ws.Cells.Item(17, 4).NumberFormat
$ws.Cells.Item(16, 5).Value = 16
$ws.Cells.Item(16, 6).Value = "Fruta"
$ws.Cells.Item(16, 7).Value = 100107
$ws.Cells.Item(16, 8).Value = "Otros"
$ws.Cells.Item(16, 9).Value = 100107002
$ws.Cells.Item(16, 10).Value = "Chirimoya"
$ws.Cells.Item(16, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(16, 12).Value = "Primera"
$ws.Cells.Item(16, 13).Value = 60
$ws.Cells.Item(16, 14).Value = 21000
$ws.Cells.Item(16, 15).Value = 21000
$ws.Cells.Item(16, 16).Value = 21000
$ws.Cells.Item(16, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(16, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(16, 19).Value = 2100
$ws.Cells.Item(16, 20).Value = 10
